$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple pairwise row swaps (columns B:AD; column A id stays fixed) ---
$row15 = $ws.Range("B15:AD15").Value2
$row16 = $ws.Range("B16:AD16").Value2
$ws.Range("B15:AD15").Value2 = $row16
$ws.Range("B16:AD16").Value2 = $row15

$row57 = $ws.Range("B57:AD57").Value2
$row58 = $ws.Range("B58:AD58").Value2
$ws.Range("B57:AD57").Value2 = $row58
$ws.Range("B58:AD58").Value2 = $row57

$row78 = $ws.Range("B78:AD78").Value2
$row79 = $ws.Range("B79:AD79").Value2
$ws.Range("B78:AD78").Value2 = $row79
$ws.Range("B79:AD79").Value2 = $row78

$row86 = $ws.Range("B86:AD86").Value2
$row87 = $ws.Range("B87:AD87").Value2
$ws.Range("B86:AD86").Value2 = $row87
$ws.Range("B87:AD87").Value2 = $row86

$row139 = $ws.Range("B139:AD139").Value2
$row140 = $ws.Range("B140:AD140").Value2
$ws.Range("B139:AD139").Value2 = $row140
$ws.Range("B140:AD140").Value2 = $row139

$row154 = $ws.Range("B154:AD154").Value2
$row155 = $ws.Range("B155:AD155").Value2
$ws.Range("B154:AD154").Value2 = $row155
$ws.Range("B155:AD155").Value2 = $row154

$row162 = $ws.Range("B162:AD162").Value2
$row163 = $ws.Range("B163:AD163").Value2
$ws.Range("B162:AD162").Value2 = $row163
$ws.Range("B163:AD163").Value2 = $row162

$row165 = $ws.Range("B165:AD165").Value2
$row166 = $ws.Range("B166:AD166").Value2
$ws.Range("B165:AD165").Value2 = $row166
$ws.Range("B166:AD166").Value2 = $row165

$row178 = $ws.Range("B178:AD178").Value2
$row180 = $ws.Range("B180:AD180").Value2
$ws.Range("B178:AD178").Value2 = $row180
$ws.Range("B180:AD180").Value2 = $row178

# --- 3-way cyclic shift among rows 177, 179, 181 ---
# after177 = before181 ; after181 = before179 ; after179 = before177
$row177 = $ws.Range("B177:AD177").Value2
$row179 = $ws.Range("B179:AD179").Value2
$row181 = $ws.Range("B181:AD181").Value2
$ws.Range("B177:AD177").Value2 = $row181
$ws.Range("B179:AD179").Value2 = $row177
$ws.Range("B181:AD181").Value2 = $row179

Write-Output "done"